# Auto-generated PowerShell Excel COM-interop script
# Applies updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.066.33"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -3.02%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.716.00"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -2.90%  "

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "310.67"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -5.68%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4785"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +6.96%  "

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3455"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.91%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "42.14"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.52%  "

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07240"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.86%  "

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.041"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -4.88%  "

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "19.75"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -4.51%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.840"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.84%  "

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "1.716.22"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -3.11%  "

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "6.816"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -5.20%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "87.48"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -5.40%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.00001033"
$cell.Style = "Normal"

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06380"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "16.45"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -3.03%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.618"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.62%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "27.112.07"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.95%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "10.75"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -4.08%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.098"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "19.94"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -1.02%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "150.38"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -6.09%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.912.54"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -3.16%  "

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.070"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.75%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "120.73"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -2.90%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.042"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -3.68%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.09219"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.602"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -2.00%  "

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.315"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -5.13%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.471"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +6.35%  "

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.02181"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -4.29%  "

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.05854"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.92%  "

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "10.92"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -7.55%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.1984"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -4.84%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.08%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "4.712"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -4.82%  "

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.5942"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -4.85%  "

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.084"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -8.00%  "

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "7.503"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -3.68%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "12.65"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -4.11%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.587"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.16%  "

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.5558"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -4.55%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "118.59"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.93%  "

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.823"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -5.73%  "

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.06645"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.91%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.086"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -4.26%  "
